$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Rename the "Employee_Details" sheet to "Employee_Personal_Details"
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Employee_Details")
$ws3.Name = "Employee_Personal_Details"

# ------------------------------------------------------------------
# 2. Update the test data on the "PIM_Add_Employee" sheet:
#    the employee used in this test case changed from
#    Sweta Arora -> Geeta Angra
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("PIM_Add_Employee")
$ws2.Range("A2").Value = "Geeta"
$ws2.Range("C2").Value = "Angra"
$ws2.Range("E2").Value = "Geeta@Angra_30"
$ws2.Range("F2").Value = "Geeta@20"
$ws2.Range("G2").Value = "Geeta@20"

# ------------------------------------------------------------------
# 3. "Employee_Personal_Details" sheet: add a new "comment" column
#    that records that personal details have been added, and tidy
#    up the header formatting (drop the now unneeded cell shading).
# ------------------------------------------------------------------

# Remove the fill shading from the "maritalstatus"/"dobmonth" headers
$ws3.Range("K1:L1").Interior.Pattern = -4142

# Remove the fill shading from the "testfield" header (force a full
# style refresh so the change actually lands on disk)
$ws3.Range("M1").Font.Bold = $false
$ws3.Range("M1").Interior.Pattern = -4142
$ws3.Range("M1").Font.Bold = $true

# New "comment" column
$ws3.Range("N1").Value = "comment"
$ws3.Range("N2").Value = "Personal Details have been added."

# Header formatting for the new column: bold text with a thin left
# border separating it from the previous column
$ws3.Range("N1").Font.Bold = $true
$ws3.Range("N1").Borders.Item(7).LineStyle = 1
$ws3.Range("N1").Borders.Item(7).Weight = 2

# Match column width to the rest of the bestFit header columns
$ws3.Range("N1").ColumnWidth = 29.2

# Move the active selection to below the newly entered data, as was
# left by the author after finishing the edit
$ws3.Range("N3").Select() | Out-Null

# Restore "PIM_Add_Employee" as the active tab (it was the active
# sheet before and after this edit)
$ws2.Activate()
